# Updates the cryptocurrency symbol list (price refresh + several rows
# reshuffled to reflect the latest coinranking.com ordering).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.30"

# Row 4
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").Value = "'3.500"
$ws.Range("E4").Value = "3LEOLEO"

# Row 5
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").Value = "'5.031"
$ws.Range("E5").Value = "4HuobiTokenHT"

# Row 6
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").Value = "'0.05615"
$ws.Range("E6").Value = "5CronosCRO"

# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'6.575"
$ws.Range("E7").Value = "6KuCoinTokenKCS"

# Row 8
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.009"
$ws.Range("E8").Value = "7GateTokenGT"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.8137"
$ws.Range("E9").Value = "8MXTokenMX"

# Row 10
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D10").Value = "'0.8360"
$ws.Range("E10").Value = "9FTXTokenFTT"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1337"
$ws.Range("E11").Value = "10WazirXWRX"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.06950"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02837"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09402"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001507"
$ws.Range("E15").Value = "14BitForexTokenBF"

# Row 16
$ws.Range("D16").Value = "'0.006146"

# Row 17
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.092"
$ws.Range("E17").Value = "16BTSETokenBTSE"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005975"
$ws.Range("E18").Value = "17OneONE"

# Row 19
$ws.Range("D19").Value = "'0.3183"

# Row 21
$ws.Range("D21").Value = "'0.1291"

# Row 22
$ws.Range("D22").Value = "'3.737"

# Row 23
$ws.Range("D23").Value = "'0.04685"

# Row 25
$ws.Range("D25").Value = "'0.001243"

# Row 26
$ws.Range("D26").Value = "'0.004531"

# Row 27
$ws.Range("D27").Value = "'0.00009691"
$ws.Range("E27").Value = "26NitroExNTX"

# Row 28
$ws.Range("D28").Value = "'0.0001938"

# Row 40
$ws.Range("D40").Value = "'0.03663"

# Row 41
$ws.Range("D41").Value = "'0.006224"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# Row 43
$ws.Range("D43").Value = "'0.002733"

# Row 44
$ws.Range("D44").Value = "'0.008173"

# Row 45
$ws.Range("D45").Value = "'0.00005295"

# Row 47
$ws.Range("D47").Value = "'0.2258"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# Row 48
$ws.Range("D48").Value = "'0.002022"
